# Update the PickupID/OrderID/CourierID style reference numbers shown in the
# "reports and screenshots" sample rows on both sheets.
#
# MXTMS!D2  (OrderID) : 11191070 -> 11191100
# MXTMS!E2  (PickupID): 7392424  -> 7392452
# Connect!AG2 (PickupID): 15595151 -> 15595207
#
# These columns store the IDs as text (not numbers), so force a text number
# format before writing the value - otherwise Excel would coerce the
# all-digit string into a numeric literal.

$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item("MXTMS")
$wsMain.Range("D2").NumberFormat = "@"
$wsMain.Range("D2").Value = "11191100"

$wsMain.Range("E2").NumberFormat = "@"
$wsMain.Range("E2").Value = "7392452"

$wsConnect = $wb.Worksheets.Item("Connect")
$wsConnect.Range("AG2").NumberFormat = "@"
$wsConnect.Range("AG2").Value = "15595207"
